$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Text columns A (Nombre) and B (NIF) for existing rows 2-4 ---
$ws.Range("A2").Value = "Wwww"
$ws.Range("A3").Value = "Yyyy"
$ws.Range("A4").Value = "Zzzz"

$ws.Range("B2").Value = "12349876W"
$ws.Range("B3").Value = "12349786Y"
$ws.Range("B4").Value = "12349687Z"

# --- Drop the old hyperlinks so they can be rebuilt cleanly ---
$ws.Hyperlinks.Delete()

# --- Rebuild hyperlinks for C2 (display text differs from the cell text) and C3 ---
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:ignacio@uniovi.es", "", "", "ww@uniovi.es")
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:nauce@uniovi.es")

# Hyperlinks.Add stamps the cell text with the display value, so fix the
# actual cell text afterwards (C3 first, then C2, per the intended order).
$ws.Range("C3").Value = "yy@uniovi.es"
$ws.Range("C2").Value = "wwuniovi.es"

# --- New row 5 (Mmmm) ---
$ws.Range("A5").Value = "Mmmm"
$ws.Range("B5").Value = "13245678M"

$ws.Hyperlinks.Add($ws.Range("C5"), "mailto:jorge@uniovi.es")
$ws.Range("C5").Value = "mm@uniovies."

# --- Numeric column D ---
$ws.Range("D2").Value = 19
$ws.Range("D3").Value = 20
$ws.Range("D4").Value = 21
$ws.Range("D5").Value = 24

# --- C4 no longer has a value or hyperlink, but keeps the hyperlink style ---
$ws.Range("C4").Value = ""

# Hyperlinks.Add() leaves behind a freshly minted cell style; reapply the
# workbook's existing "Hipervínculo" style so these cells match the others.
$ws.Range("C2").Style = "Hipervínculo"
$ws.Range("C3").Style = "Hipervínculo"
$ws.Range("C4").Style = "Hipervínculo"
$ws.Range("C5").Style = "Hipervínculo"

# --- Selection moves to C3 ---
$null = $ws.Range("C3").Select()
